$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.747.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.776.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.212.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.769.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.936"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.700.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0971"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0451"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.83%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.05%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("E45").Value = "  -7.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.084.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("E51").Value = "  +1.40%  "
